# feat: add 2022-Q1 data
#
# 1) Insert a new "2022-Q1" sheet (fund-holding detail) right before the
#    "总计" (summary) sheet.
# 2) Insert a new summary row ("2022-Q1", 6, 1.06) at the top of the
#    "总计" sheet's data, shifting the existing quarters down.

$wb = $excel.ActiveWorkbook
$totalSheetAnchor = $wb.Worksheets.Item("总计")
$q4Sheet = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 1) New "2022-Q1" detail sheet, positioned right before "总计"
# ---------------------------------------------------------------------
# NOTE: Worksheets.Add(beforeSheet) re-seats the engine's sheet anchors by
# position, so $totalSheetAnchor itself ends up referring to the freshly
# inserted sheet once it occupies that slot. Do all the new-sheet setup
# through $ws, then look "总计" back up by name afterwards (see below)
# rather than continuing to use $totalSheetAnchor.
$ws = $wb.Worksheets.Add($totalSheetAnchor)
$ws.Name = "2022-Q1"
$ws.Outline.SummaryRow = 1
$ws.Outline.SummaryColumn = 1

# Clone the look (bold + border header row, bold + border index column)
# from the existing "2021-Q4" sheet, which has the same 8-column layout.
$q4Sheet.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$q4Sheet.Range("A2:H5").Copy()
$ws.Range("A2:H7").PasteSpecial(-4122)

# Header row
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Data rows: idx, fund code, fund name, fund size, stock position,
# position ratio, holding market value (100M yuan), position rank.
$rows = @(
    @(0, "007110", "国投瑞银港股通价值发现混合",         "23.33", "93.33", "2.32", "0.5413", 10),
    @(1, "010010", "国投瑞银港股通6个月定期开放股票",     "8.09",  "93.58", "3.15", "0.2548", 8),
    @(2, "011157", "弘毅远方港股通智选领航混合A",         "3.65",  "90.35", "4.18", "0.1526", 8),
    @(3, "005646", "中海沪港深多策略灵活配置混合",         "1.78",  "88.15", "4.30", "0.0765", 5),
    @(4, "011158", "弘毅远方港股通智选领航混合C",         "0.78",  "90.35", "4.18", "0.0326", 8),
    @(5, "005770", "信达澳银中证沪港深高股息精选指数",     "0.01",  "92.47", "2.57", "0.0003", 5)
)

$r = 2
foreach ($row in $rows) {
    $ws.Range("A$r").Value = $row[0]
    # The leading apostrophe forces these numeric-looking values to be
    # stored as text (matching the target's t="inlineStr"/shared-string
    # cells for the fund size / position / market-value columns).
    $ws.Range("B$r").Value = "'" + $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = "'" + $row[3]
    $ws.Range("E$r").Value = "'" + $row[4]
    $ws.Range("F$r").Value = "'" + $row[5]
    $ws.Range("G$r").Value = "'" + $row[6]
    $ws.Range("H$r").Value = $row[7]

    # Entering text via the leading apostrophe auto-applies a "Text"
    # number-format style to that cell. The source sheet leaves these
    # data cells style-less (only the header row and column A carry
    # style "2"), so re-paste formats from the always style-less name
    # column (C) onto the rest of the row to strip that back off.
    $ws.Range("C$r").Copy()
    $ws.Range("B$r").PasteSpecial(-4122)
    $ws.Range("D$r").PasteSpecial(-4122)
    $ws.Range("E$r").PasteSpecial(-4122)
    $ws.Range("F$r").PasteSpecial(-4122)
    $ws.Range("G$r").PasteSpecial(-4122)

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) "总计" sheet: add a "2022-Q1" row at the top of the data, shifting
#    the previously existing quarters down by one row.
# ---------------------------------------------------------------------
# Re-fetch "总计" by name now that the sheet collection has shifted
# (see note above) so this reliably points at the real summary sheet.
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Range("A2:D2").Insert()

# Re-apply the index-column style (bold + border) that Insert() doesn't
# carry onto the new row, and clear the stray style Insert() leaves on
# the rest of the new row so it matches the other data rows.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 6
$totalSheet.Range("D2").Value = 1.06

# Renumber the shifted-down rows' index column (0-based, sequential).
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4

Write-Output "2022-Q1 sheet added; 总计 updated"
